$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the PATENT NUMBER column (B2:B9): each list has its last
# (trailing) patent number removed, i.e. the item-based user data is
# reduced by one entry per row.
$ws.Range("B2").Value = "10349422, 10349418"
$ws.Range("B3").Value = "10346095, 10346094"
$ws.Range("B4").Value = "9967960, 9967646"
$ws.Range("B5").Value = "9967277, 9965766"
$ws.Range("B6").Value = "10331583, 10327202"
$ws.Range("B7").Value = "9794808, 9794797"
$ws.Range("B8").Value = "10312751, 10312750"
$ws.Range("B9").Value = "9967277, 9965766"

# Move the active selection to B5 as in the saved workbook.
$ws.Range("B5").Select()
